$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "LAST SCRAPE DATE" column (F) for every data row (2-83):
# the whole column shares one string value, 2019-03-07 -> 2019-03-12.
# Force the column to text first so Excel keeps the value as a plain
# string (matching the source data) instead of reinterpreting it as a date.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row
$dateRange = $ws.Range("F2:F$lastRow")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2019-03-12"

# The remaining changes are per-row updates to GAME NAME (C), and the two
# numeric columns D (total top prizes) / E (top prizes remaining).

# Row 4: $1 MONOPOLY JACKPOT - only remaining count changes
$ws.Range("E4").Value = 68

# Rows 5 & 6 swap their data (TRIPLE WIN <-> HEADS OR TAILS)
$ws.Range("C5").Value = "HEADS OR TAILS"
$ws.Range("D5").Value = 1395
$ws.Range("E5").Value = 37

$ws.Range("C6").Value = "TRIPLE WIN"
$ws.Range("D6").Value = 1360
$ws.Range("E6").Value = 30

# Rows 9 & 10 swap their data (MONEY BAGS <-> $100 IN A FLASH)
$ws.Range("C9").Value = "$100 IN A FLASH"
$ws.Range("D9").Value = 1380
$ws.Range("E9").Value = 13

$ws.Range("C10").Value = "MONEY BAGS"
$ws.Range("D10").Value = 1412
$ws.Range("E10").Value = 372

# Row 11: $10,000 GOLD RUSH CLASSIC - only remaining count changes
$ws.Range("E11").Value = 71

# Rows 27, 28, 29 rotate (TRIPLE CA$H -> $50,000 FLAMINGO MULTIPLIER ->
# $50,000 GOLD RUSH CLASSIC -> TRIPLE CA$H)
$ws.Range("C27").Value = "$50,000 FLAMINGO MULTIPLIER"
$ws.Range("D27").Value = 1297
$ws.Range("E27").Value = 13

$ws.Range("C28").Value = "$50,000 GOLD RUSH CLASSIC"
$ws.Range("D28").Value = 1409
$ws.Range("E28").Value = 57

$ws.Range("C29").Value = 'TRIPLE CA$H'
$ws.Range("D29").Value = 1413
$ws.Range("E29").Value = 9

# Row 43: SCRABBLE - only remaining count changes
$ws.Range("E43").Value = 11

# Rows 55, 56, 57 rotate ($1,000,000 GOLD RUSH CLASSIC -> BONUS DOUBLE MATCH
# -> WINNING STREAK -> $1,000,000 GOLD RUSH CLASSIC)
$ws.Range("C55").Value = "BONUS DOUBLE MATCH"
$ws.Range("D55").Value = 1348
$ws.Range("E55").Value = 2

$ws.Range("C56").Value = "WINNING STREAK"
$ws.Range("D56").Value = 1302
$ws.Range("E56").Value = 7

$ws.Range("C57").Value = "$1,000,000 GOLD RUSH CLASSIC"
$ws.Range("D57").Value = 1410
$ws.Range("E57").Value = 26
